$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 9002
$ws.Range("B6").Value = "BodyWeapon"
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.1
$ws.Range("I6").Value = 0.1

$ws.Range("J6").Select()
